$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''245.89'
$ws.Range('E2').Value = '''-0.36%'

$ws.Range('D3').Value = '''29.79'
$ws.Range('E3').Value = '''-1.05%'

$ws.Range('D4').Value = '''5.155'
$ws.Range('E4').Value = '''-0.18%'

$ws.Range('D5').Value = '''0.05774'
$ws.Range('E5').Value = '''0.75%'

$ws.Range('D6').Value = '''6.652'
$ws.Range('E6').Value = '''0.98%'

$ws.Range('D7').Value = '''3.234'
$ws.Range('E7').Value = '''6.50%'

$ws.Range('D8').Value = '''0.8493'
$ws.Range('E8').Value = '''-0.87%'

$ws.Range('D9').Value = '''0.8532'
$ws.Range('E9').Value = '''-2.72%'

$ws.Range('D10').Value = '''0.1375'
$ws.Range('E10').Value = '''0.54%'

$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = '''0.03375'
$ws.Range('E11').Value = '''1.65%'

$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.07086'
$ws.Range('E12').Value = '''1.30%'

$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03252'
$ws.Range('E13').Value = '''11.09%'

$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09370'
$ws.Range('E14').Value = '''-0.09%'

$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001531'
$ws.Range('E15').Value = '''-0.05%'

$ws.Range('B16').Value = 'One'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D16').Value = '''0.0005982'
$ws.Range('E16').Value = '''-0.40%'

$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').Value = '''0.006035'
$ws.Range('E17').Value = '''-0.25%'

$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').Value = '''3.509'
$ws.Range('E18').Value = '''0.00%'

$ws.Range('B19').Value = 'BTSEToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D19').Value = '''2.221'
$ws.Range('E19').Value = '''1.75%'

$ws.Range('B20').Value = 'BitpandaEcosystemToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D20').Value = '''0.3159'
$ws.Range('E20').Value = '''0.48%'

$ws.Range('D21').Value = '''0.1300'
$ws.Range('E21').Value = '''-0.36%'

$ws.Range('D22').Value = '''3.497'
$ws.Range('E22').Value = '''-2.77%'

$ws.Range('B23').Value = 'ZBToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D23').Value = '''0.1410'
$ws.Range('E23').Value = '''2.52%'

$ws.Range('B24').Value = 'CoinExToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D24').Value = '''0.04126'
$ws.Range('E24').Value = '''-0.85%'

$ws.Range('D25').Value = '''0.001227'
$ws.Range('E25').Value = '''1.13%'

$ws.Range('D26').Value = '''0.004141'
$ws.Range('E26').Value = '''-7.98%'

$ws.Range('E27').Value = '''1.94%'

$ws.Range('D40').Value = '''0.03755'
$ws.Range('E40').Value = '''-0.74%'

$ws.Range('D41').Value = '''0.1070'
$ws.Range('E41').Value = '''0.17%'

$ws.Range('D42').Value = '''0.002471'
$ws.Range('E42').Value = '''-2.14%'

$ws.Range('D43').Value = '''0.003525'
$ws.Range('E43').Value = '''-38.80%'

$ws.Range('D44').Value = '''0.009511'
$ws.Range('E44').Value = '''-4.83%'

$ws.Range('D45').Value = '''0.00005428'
$ws.Range('E45').Value = '''6.85%'

$ws.Range('E46').Value = '''0.28%'

$ws.Range('E47').Value = '''-20.01%'

$ws.Range('D48').Value = '''0.002196'
$ws.Range('E48').Value = '''-19.27%'

$ws.Range('E49').Value = '''0.28%'

$ws.Range('E50').Value = '''0.28%'
